# Fix - Fixed birthday input in the spreadsheet.
#
# The "Fecha_de_Nacimiento" (birthdate) column D had some entries stored as
# real Excel dates (serial numbers) while others were plain text strings
# like "22/12/1976". This mix caused inconsistent display. The fix re-enters
# every birthdate in column D as literal dd/mm/yyyy TEXT (number format
# "Text"), including zero-padding the day on the one entry that was missing
# it ("14/9/1977" -> "14/09/1977").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make the whole header+data range of column D a Text-formatted column so
# every birthdate (old and new) renders as literal text, not a date serial.
$ws.Range("D1:D8").NumberFormat = "@"

# Re-type the birthdates that used to be stored as real dates, as plain
# dd/mm/yyyy text.
$ws.Range("D8").Value = "01/09/2017"
$ws.Range("D4").Value = "20/02/1990"
$ws.Range("D2").Value = "24/12/1999"

# Normalize the one text birthdate that was missing zero-padding on the day.
$ws.Range("D3").Value = "14/09/1977"

# Column C ("Nombre") is now wide enough that Excel remembers an explicit
# best-fit width for it.
$ws.Columns("C").AutoFit()

# Tidy up stray formatting left behind below the data (an empty
# date-formatted cell in D10 is cleared, and the empty right-aligned
# placeholder cell that used to sit at F13 is now at D13 instead).
$ws.Range("D10").Clear()
$ws.Range("F13").Clear()

$ws.Range("D13").NumberFormat = "mm-dd-yy"
$ws.Range("D13").HorizontalAlignment = -4152

# A handful of other empty cells further down picked up the underlined
# "hyperlink-like" formatting while the user was clicking around checking
# things, leaving empty formatted placeholders behind.
$ws.Range("F16").Font.Underline = $true
$ws.Range("F19").Font.Underline = $true
$ws.Range("H22").Font.Underline = $true
$ws.Range("H23").Font.Underline = $true

$ws.Range("H22:H23").Select()
